$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2980.5
$ws.Range("J17").Value = 3255.625
$ws.Range("L17").Value = 9766.875
$ws.Range("N17").Value = -10102.875
$ws.Range("H33").Value = 369.26315
$ws.Range("I33").Value = 373.1111
$ws.Range("K33").Value = 373.1111
$ws.Range("M33").Value = -144.1111
$ws.Range("H107").Value = 713.3333
$ws.Range("I107").Value = 713.3333
$ws.Range("K107").Value = 713.3333
$ws.Range("M107").Value = 1206.6667
$ws.Range("H116").Value = 10000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 10000
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -16884
$ws.Range("H132").Value = 3130.868
$ws.Range("I132").Value = 3154.5386
$ws.Range("K132").Value = 9463.6158
$ws.Range("M132").Value = -6933.6158
$ws.Range("H137").Value = 4859.091
$ws.Range("I137").Value = 3955.8
$ws.Range("J137").Value = 6794.7144
$ws.Range("K137").Value = 11867.4
$ws.Range("L137").Value = 20384.1432
$ws.Range("M137").Value = -9317.400000000001
$ws.Range("N137").Value = -25484.1432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 31696.215
$ws.Range("I45").Value = 43591.9
$ws.Range("K45").Value = 43591.9
$ws.Range("M45").Value = -43214.9
$ws.Range("H74").Value = 279634.2
$ws.Range("I74").Value = 398341.84
$ws.Range("J74").Value = 2649.6667
$ws.Range("K74").Value = 398341.84
$ws.Range("L74").Value = 2649.6667
$ws.Range("M74").Value = -397467.84
$ws.Range("N74").Value = -4397.6667
$ws.Range("H76").Value = 100000
$ws.Range("J76").Value = 100000
$ws.Range("L76").Value = 100000
$ws.Range("N76").Value = -100676
$ws.Range("H77").Value = 279634.2
$ws.Range("I77").Value = 398341.84
$ws.Range("J77").Value = 2649.6667
$ws.Range("K77").Value = 1991709.2
$ws.Range("L77").Value = 13248.3335
$ws.Range("M77").Value = -1987341.2
$ws.Range("N77").Value = -21984.3335
$ws.Range("H79").Value = 100000
$ws.Range("J79").Value = 100000
$ws.Range("L79").Value = 100000
$ws.Range("N79").Value = -102340
$ws.Range("H110").Value = 801.125
$ws.Range("I110").Value = 684
$ws.Range("J110").Value = 996.3333
$ws.Range("K110").Value = 684
$ws.Range("L110").Value = 996.3333
$ws.Range("M110").Value = 1361
$ws.Range("N110").Value = -5086.3333
$ws.Range("H122").Value = 3330.1052
$ws.Range("J122").Value = 3450.0908
$ws.Range("L122").Value = 10350.2724
$ws.Range("N122").Value = -15250.2724
$ws.Range("H132").Value = 3056.5789
$ws.Range("J132").Value = 2748.7
$ws.Range("L132").Value = 8246.099999999999
$ws.Range("N132").Value = -13306.1

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 65002624
$ws.Range("I105").Value = 10000000
$ws.Range("J105").Value = 83336830
$ws.Range("K105").Value = 10000000
$ws.Range("L105").Value = 83336830
$ws.Range("M105").Value = -9998253
$ws.Range("N105").Value = -83340324
$ws.Range("H107").Value = 2080061.4
$ws.Range("I107").Value = 2565074.5
$ws.Range("J107").Value = 1434.1428
$ws.Range("K107").Value = 2565074.5
$ws.Range("L107").Value = 1434.1428
$ws.Range("M107").Value = -2563154.5
$ws.Range("N107").Value = -5274.1428
$ws.Range("H139").Value = 101606.336
$ws.Range("J139").Value = 101606.336
$ws.Range("L139").Value = 101606.336
$ws.Range("N139").Value = -111886.336

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4145.3193
$ws.Range("I31").Value = 4854.846
$ws.Range("J31").Value = 3874.0293
$ws.Range("K31").Value = 4854.846
$ws.Range("L31").Value = 3874.0293
$ws.Range("M31").Value = -4559.846
$ws.Range("N31").Value = -4464.0293
$ws.Range("H34").Value = 4145.3193
$ws.Range("I34").Value = 4854.846
$ws.Range("J34").Value = 3874.0293
$ws.Range("K34").Value = 4854.846
$ws.Range("L34").Value = 3874.0293
$ws.Range("M34").Value = -4652.846
$ws.Range("N34").Value = -4278.0293
$ws.Range("H51").Value = 59999.5
$ws.Range("J51").Value = 59999.5
$ws.Range("L51").Value = 59999.5
$ws.Range("N51").Value = -61471.5
$ws.Range("H61").Value = 59999.5
$ws.Range("J61").Value = 59999.5
$ws.Range("L61").Value = 59999.5
$ws.Range("N61").Value = -60695.5
$ws.Range("H107").Value = 2632408.2
$ws.Range("I107").Value = 4167305.8
$ws.Range("J107").Value = 1155
$ws.Range("K107").Value = 4167305.8
$ws.Range("L107").Value = 1155
$ws.Range("M107").Value = -4165385.8
$ws.Range("N107").Value = -4995
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H134").Value = 2697.9756
$ws.Range("I134").Value = 2363.5278
$ws.Range("K134").Value = 7090.5834
$ws.Range("M134").Value = -4555.5834

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 2609.889
$ws.Range("I26").Value = 3619.8333
$ws.Range("J26").Value = 590
$ws.Range("K26").Value = 10859.4999
$ws.Range("L26").Value = 1770
$ws.Range("M26").Value = -10571.4999
$ws.Range("N26").Value = -2346
$ws.Range("H46").Value = 403324.38
$ws.Range("J46").Value = 734.8
$ws.Range("L46").Value = 2204.4
$ws.Range("N46").Value = -2386.4
$ws.Range("H75").Value = 1500
$ws.Range("J75").Value = 1500
$ws.Range("L75").Value = 4500
$ws.Range("N75").Value = -6496
$ws.Range("H78").Value = 1500
$ws.Range("J78").Value = 1500
$ws.Range("L78").Value = 13500
$ws.Range("N78").Value = -23484
$ws.Range("H86").Value = 2992.6667
$ws.Range("J86").Value = 2935.111
$ws.Range("L86").Value = 8805.332999999999
$ws.Range("N86").Value = -11177.333
$ws.Range("H89").Value = 2992.6667
$ws.Range("J89").Value = 2935.111
$ws.Range("L89").Value = 26415.999
$ws.Range("N89").Value = -38271.999
$ws.Range("H134").Value = 2665.25
$ws.Range("I134").Value = 1664.6428
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 4993.928400000001
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = 76.07159999999931
$ws.Range("N134").Value = -25140
$ws.Range("H140").Value = 3311.0356
$ws.Range("I140").Value = 3100.7896
$ws.Range("J140").Value = 3754.889
$ws.Range("K140").Value = 9302.3688
$ws.Range("L140").Value = 11264.667
$ws.Range("M140").Value = -4122.3688
$ws.Range("N140").Value = -21624.667

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H58").Value = 49999.5
$ws.Range("I58").Value = 49999.5
$ws.Range("K58").Value = 49999.5
$ws.Range("M58").Value = -49722.5
$ws.Range("H97").Value = 1817.8148
$ws.Range("I97").Value = 1355.7391
$ws.Range("J97").Value = 4474.75
$ws.Range("K97").Value = 1355.7391
$ws.Range("L97").Value = 4474.75
$ws.Range("M97").Value = -859.7391
$ws.Range("N97").Value = -5466.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1779.9
$ws.Range("I46").Value = 1366.5555
$ws.Range("K46").Value = 1366.5555
$ws.Range("M46").Value = -1178.5555
$ws.Range("H68").Value = 4437.3076
$ws.Range("I68").Value = 3398.5
$ws.Range("J68").Value = 6099.4
$ws.Range("K68").Value = 3398.5
$ws.Range("L68").Value = 6099.4
$ws.Range("M68").Value = -2649.5
$ws.Range("N68").Value = -7597.4
$ws.Range("H71").Value = 4437.3076
$ws.Range("I71").Value = 3398.5
$ws.Range("J71").Value = 6099.4
$ws.Range("K71").Value = 16992.5
$ws.Range("L71").Value = 30497
$ws.Range("M71").Value = -13248.5
$ws.Range("N71").Value = -37985
$ws.Range("H122").Value = 11532.267
$ws.Range("I122").Value = 7748.75
$ws.Range("K122").Value = 23246.25
$ws.Range("M122").Value = -20796.25
$ws.Range("H132").Value = 4581.391
$ws.Range("I132").Value = 3942.6
$ws.Range("K132").Value = 11827.8
$ws.Range("M132").Value = -9297.799999999999
$ws.Range("H136").Value = 2155.75
$ws.Range("I136").Value = 1277
$ws.Range("K136").Value = 3831
$ws.Range("M136").Value = -1281

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 14999.5
$ws.Range("I52").Value = 14999.5
$ws.Range("K52").Value = 14999.5
$ws.Range("M52").Value = -14773.5
$ws.Range("H62").Value = 6887.6665
$ws.Range("I62").Value = 6048.75
$ws.Range("J62").Value = 7558.8
$ws.Range("K62").Value = 6048.75
$ws.Range("L62").Value = 7558.8
$ws.Range("M62").Value = -5424.75
$ws.Range("N62").Value = -8806.799999999999
$ws.Range("H65").Value = 6887.6665
$ws.Range("I65").Value = 6048.75
$ws.Range("J65").Value = 7558.8
$ws.Range("K65").Value = 30243.75
$ws.Range("L65").Value = 37794
$ws.Range("M65").Value = -27123.75
$ws.Range("N65").Value = -44034
$ws.Range("H115").Value = 49999
$ws.Range("J115").Value = 49999
$ws.Range("L115").Value = 49999
$ws.Range("N115").Value = -53133
$ws.Range("H138").Value = 84527
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 84527
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -94807
